# Added Configurable zero_before_threshold parameter to enable setting dims
# before noise_threshold or First Rise Point to 0.
#
# This updates the cached First_Noticeable_Increase_Index (C),
# First_Noticeable_Increase_Cumulative_Value (E) and Pulse_Width (G)
# columns on each Step3_DataPts_* sheet to reflect the new calculation
# results produced after that feature was introduced.

$wb = $excel.ActiveWorkbook

# Values are identical across all four Step3_DataPts_* sheets for columns
# C (First_Noticeable_Increase_Index) and E (First_Noticeable_Increase_Cumulative_Value),
# since those do not depend on the Intensity_Threshold. Column G (Pulse_Width)
# depends on both the new First_Noticeable_Increase_Index and the (unchanged)
# Point_Exceeds_Index, so it differs by sheet/threshold.

$updates = @{
    "Step3_DataPts_0.5" = @(
        @{ Row = 2; C = 90; E = 0.04615289647583285; G = 34 },
        @{ Row = 3; C = 91; E = 0.03959320895107479; G = 52 },
        @{ Row = 4; C = 87; E = 0.06818682133623107; G = 44 },
        @{ Row = 5; C = 87; E = 0.02313440193116293; G = 44 },
        @{ Row = 6; C = 87; E = 0.04659958260154852; G = 36 }
    )
    "Step3_DataPts_0.7" = @(
        @{ Row = 2; C = 90; E = 0.04615289647583285; G = 57 },
        @{ Row = 3; C = 91; E = 0.03959320895107479; G = 70 },
        @{ Row = 4; C = 87; E = 0.06818682133623107; G = 64 },
        @{ Row = 5; C = 87; E = 0.02313440193116293; G = 63 },
        @{ Row = 6; C = 87; E = 0.04659958260154852; G = 56 }
    )
    "Step3_DataPts_0.8" = @(
        @{ Row = 2; C = 90; E = 0.04615289647583285; G = 69 },
        @{ Row = 3; C = 91; E = 0.03959320895107479; G = 78 },
        @{ Row = 4; C = 87; E = 0.06818682133623107; G = 71 },
        @{ Row = 5; C = 87; E = 0.02313440193116293; G = 70 },
        @{ Row = 6; C = 87; E = 0.04659958260154852; G = 71 }
    )
    "Step3_DataPts_0.9" = @(
        @{ Row = 2; C = 90; E = 0.04615289647583285; G = 93 },
        @{ Row = 3; C = 91; E = 0.03959320895107479; G = 98 },
        @{ Row = 4; C = 87; E = 0.06818682133623107; G = 96 },
        @{ Row = 5; C = 87; E = 0.02313440193116293; G = 94 },
        @{ Row = 6; C = 87; E = 0.04659958260154852; G = 96 }
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($rowUpdate in $updates[$sheetName]) {
        $r = $rowUpdate.Row
        $ws.Cells.Item($r, 3).Value = $rowUpdate.C   # Column C
        $ws.Cells.Item($r, 5).Value = $rowUpdate.E   # Column E
        $ws.Cells.Item($r, 7).Value = $rowUpdate.G   # Column G
    }
}

$wb.Save()
